$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header text, written in the same order the author authored them in
# (Email, then the two University Number columns) so the shared-string
# table grows in that order.
$ws.Range("I1").Value = "Email"
$ws.Range("A1").Value = "University Roll Number"
$ws.Range("B1").Value = "University Registration Number"

# Match the header row styling (bold, size 14 -- same as the rest of row 1)
$ws.Range("A1").Style = $ws.Range("C1").Style
$ws.Range("I1").Style = $ws.Range("C1").Style

# Column widths
$ws.Columns.Item(1).ColumnWidth = 27.6640625
$ws.Columns.Item(2).ColumnWidth = 37.21875
$ws.Columns.Item(9).ColumnWidth = 15

# Clear the stale selection left on the sheet view
$ws.Range("A1").Select()
